$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Target change (paragraph "ตาราง … Activity Diagram"):
#   "ตาราง" + " "   + "… "  + "Activity" + " Diagram"
#   -> "ตาราง" + "ที่ " + "1" + " "       + "Activity" + " Diagram"
#
# i.e. the run holding the single space right after "ตาราง" becomes
# "ที่ ", and the run holding the ellipsis placeholder "… " is split
# into a "1" run and a plain " " run (matching the diff, which shows
# the ellipsis run text changing to "1" and a brand-new run carrying
# the trailing space).
# ------------------------------------------------------------------

# Locate the "ตาราง" run and the space run that currently sits right
# after it, plus the "… " run right after that (three separate,
# adjacent runs).
$rWord = $d.Content
$rWord.Find.Execute("ตาราง", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$wordStart = $rWord.Start
$wordEnd = $rWord.End

$rEllipsis = $d.Content
$rEllipsis.Find.Execute("… ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$ellipsisStart = $rEllipsis.Start
$ellipsisEnd = $rEllipsis.End

# --- Edit 1: replace the lone space after "ตาราง" with "ที่ " -------
# Use a transient (no-op) bold toggle on "ตาราง" itself so that, once
# the replacement text is typed, Word does not silently fold the
# preceding "ตาราง" run into the run we are about to rewrite (they
# currently share identical run formatting). Scope the toggle tightly
# to the "ตาราง" run only so unrelated text earlier in the document is
# left untouched.
$rProtectWord = $d.Range($wordStart, $wordEnd)
$rProtectWord.Font.Bold = $true
$rProtectWord.Font.Bold = $false

$rSpace = $d.Range($wordEnd, $ellipsisStart)
$rSpace.Text = "ที่ "

$rProtectWord2 = $d.Range($wordStart, $wordEnd)
$rProtectWord2.Font.Bold = $true
$rProtectWord2.Font.Bold = $false

# --- Edit 2: turn "… " into "1" + " " -------------------------------
$rEllipsis2 = $d.Content
$rEllipsis2.Find.Execute("… ", $true, $false, $false, $false, $false, $true, 1, $false, "1 ", 2) | Out-Null

# After the two text edits above, Word's run model may have coalesced
# several neighbouring runs that share identical formatting (e.g. "1",
# " ", "Activity", " Diagram" all use the same rPr). Re-split them back
# into individual runs - matching the target markup - using harmless
# bold on/off toggles (Word leaves no trace in the saved XML when a
# property is toggled back to its original value).
$rOne = $d.Content
$rOne.Find.Execute("1 Activity Diagram", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$segStart = $rOne.Start
$segEnd = $rOne.End

$oneEnd = $segStart + 1          # end of "1"
$spaceEnd = $oneEnd + 1          # end of the following " "
$activityEnd = $spaceEnd + 8     # end of "Activity"

$splitRanges = @(
    @($segStart, $oneEnd),
    @($oneEnd, $spaceEnd),
    @($spaceEnd, $activityEnd),
    @($activityEnd, $segEnd)
)
foreach ($bounds in $splitRanges) {
    $rSplit = $d.Range($bounds[0], $bounds[1])
    $rSplit.Font.Bold = $true
    $rSplit.Font.Bold = $false
}
